$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Devika "
$ws.Range("B3").Value = "Ugle"
$ws.Range("A1").Value = "Darshan"
$ws.Range("B1").Value = "Ghurde"

$ws.Range("B3").Font.Bold = $true

$ws.Range("E3").Select()
